$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("ListOfValues")

# --- ListOfValues: add the Bybit_Testnet exchange entry ---
$ws2.Range("A4").Value = "Bybit_Testnet"

# --- Sheet1: fix UltimateScalper test row, remove stale duplicate rows ---

# Remove the three duplicate/obsolete backtest rows (rows 3:5) - everything
# below shifts up to fill the gap.
$ws1.Range("A3:A5").EntireRow.Delete()

# Update the "From" date for the remaining UltimateScalper test row.
$ws1.Range("D2").Value = 44593

# Set the optional strategy-settings JSON for the test (previously blank).
$ws1.Range("K2").Value = '{"EMA_Fast": 9, "EMA_Slow": 55, "EMA_Trend": 200, "RSI": 4, "RSI_Low": 19, "RSI_High": 81, "ADX": 17, "ADX_Threshold": 24, "MACD_Fast": 12, "MACD_Slow": 24, "MACD_Signal": 9, "BB_Length": 34, "BB_Mult": 1}'

# --- View / selection bookkeeping, matches the saved workbook state ---
$ws2.Range("B23").Select()

$ws3 = $wb.Worksheets.Item("StrategyDictionaries")
$ws3.Activate()
$ws3.Range("B4").Select()
$ws3.Application.ActiveWindow.ScrollColumn = 2

$ws1.Activate()
$ws1.Range("D3").Select()
$ws1.Application.ActiveWindow.ScrollColumn = 2
